$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in newly entered scores that were previously blank
$ws.Range("G7").Value = 20
$ws.Range("L7").Value = 15
$ws.Range("B18").Value = 10

# Move the active selection to B19 (reflecting where the user clicked next)
$ws.Range("B19").Select()
